$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 2.35
$ws.Range("C3").Value = 4.75
$ws.Range("D3").Value = 0.35
$ws.Range("E3").Value = 1.1499999999999999

$ws.Range("H3").Value = 2.4
$ws.Range("I3").Value = 4.8
$ws.Range("J3").Value = 0.4
$ws.Range("K3").Value = 1.2

$ws.Range("B4").Value = 1.1499999999999999
$ws.Range("C4").Value = 2.85
$ws.Range("D4").Value = 0.15
$ws.Range("E4").Value = 0.55000000000000004

$ws.Range("H4").Value = 1.2
$ws.Range("I4").Value = 2.9
$ws.Range("J4").Value = 0.2
$ws.Range("K4").Value = 0.6

$ws.Range("B5").Value = 0.55000000000000004
$ws.Range("C5").Value = 1.1499999999999999
$ws.Range("D5").Value = 0.06
$ws.Range("E5").Value = 0.25

$ws.Range("H5").Value = 0.6
$ws.Range("I5").Value = 1.2
$ws.Range("J5").Value = 0.1
$ws.Range("K5").Value = 0.3

$ws.Range("E5").Select()
